# Weekly update: a new price observation (week) is inserted as the new
# row 181 on the "Papa" sheet; every row that used to be 181-185 shifts
# down by one (to 182-186), keeping its data unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 181, pushing existing rows 181-185 down to 182-186.
$ws.Rows.Item(181).Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Cells.Item(181, 1).Value  = 7
$ws.Cells.Item(181, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(181, 3).Value  = "Ñuble"
$ws.Cells.Item(181, 4).Value  = 44448
$ws.Cells.Item(181, 5).Value  = 16
$ws.Cells.Item(181, 6).Value  = 100114001
$ws.Cells.Item(181, 7).Value  = "Papa"
$ws.Cells.Item(181, 8).Value  = "Patagonia"
$ws.Cells.Item(181, 9).Value  = "1a (guarda)"
$ws.Cells.Item(181, 10).Value = 300
$ws.Cells.Item(181, 11).Value = 7000
$ws.Cells.Item(181, 12).Value = 7500
$ws.Cells.Item(181, 13).Value = 7250
$ws.Cells.Item(181, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(181, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(181, 16).Value = 290
$ws.Cells.Item(181, 17).Value = 25
$ws.Cells.Item(181, 18).Value = "Hortaliza"
